# Update the customer sample data (rows 2-4) with a new set of sample
# registrants, and trim the Address column so the wrapped text again fits
# the (now shorter) row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Maureen J. Delafuente
$ws.Range("B2").Value = "Maureen"
$ws.Range("C2").Value = "J."
$ws.Range("D2").Value = "MaureenJDelafuente@jourrapide.com"
$ws.Range("E2").Value = "971 Alpha AvenueMarshall, TX 75670"
$ws.Range("F2").Value = "903-938-1180"

# Row 3 - Arlene K. Simon
$ws.Range("B3").Value = "Arlene"
$ws.Range("C3").Value = "K."
$ws.Range("D3").Value = "ArleneKSimon@jourrapide.com"
$ws.Range("E3").Value = "4214 Reynolds AlleyBellflower, CA 90706"
$ws.Range("F3").Value = "562-285-4192"

# Row 4 - Judy D. Allen
$ws.Range("B4").Value = "Judy"
$ws.Range("C4").Value = "D."
$ws.Range("D4").Value = "JudyDAllen@jourrapide.com"
$ws.Range("E4").Value = "1403 Half and Half DriveKennedy Meadows, CA 93527"
$ws.Range("F4").Value = "559-850-7665"

# Row heights re-flow now that the (wrapped) Address text changed length.
$ws.Rows.Item(2).RowHeight = 72.5
$ws.Rows.Item(3).RowHeight = 72.5
$ws.Rows.Item(4).RowHeight = 101.5
